# Rename the existing sheet and add a second scenario sheet, mirroring the
# "Add two numbers" example but with negative operands.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "SampleTest1" -> "Test1" ------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Test1"
[void]$ws1.Range("G43").Select()

# --- Sheet 2: new "Test2" sheet, inserted right after Test1 ---------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Test2"

$ws2.Range("A1").Value = "Scenario"
$ws2.Range("B1").Value = "Add two negative numbers"

$ws2.Range("A2").Value = "Given"
$ws2.Range("B2").Value = "the first number is"
$ws2.Range("C2").Value = -50

$ws2.Range("A3").Value = "And"
$ws2.Range("B3").Value = "the second number is"
$ws2.Range("C3").Value = -70

$ws2.Range("A4").Value = "When"
$ws2.Range("B4").Value = "the two numbers are added"

$ws2.Range("A5").Value = "Then"
$ws2.Range("B5").Value = "the result should be"
$ws2.Range("C5").Formula = "=C2+C3"

[void]$ws2.Range("A1:C5").Select()
